$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.707.55'
$ws.Range('E2').Value = '  +4.83%  '
$ws.Range('D3').Value = '2.978.93'
$ws.Range('E3').Value = '  +2.84%  '
$ws.Range('E4').Value = '  +0.23%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '579.25'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('E5').Value = '  +1.99%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '152.15'
$ws.Range('D6').NumberFormat = 'General'
$ws.Range('E6').Value = '  +6.56%  '
$ws.Range('E7').Value = '  -0.14%  '
$ws.Range('D8').Value = '2.975.58'
$ws.Range('E8').Value = '  +2.87%  '
$ws.Range('E9').Value = '  +0.37%  '
$ws.Range('E10').Value = '  +5.35%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.152'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('E11').Value = '  +2.93%  '
$ws.Range('E12').Value = '  +2.51%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000238'
$ws.Range('D13').NumberFormat = 'General'
$ws.Range('E13').Value = '  +2.64%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.16'
$ws.Range('D14').NumberFormat = 'General'
$ws.Range('E14').Value = '  +6.48%  '
$ws.Range('E15').Value = '  +0.73%  '
$ws.Range('D16').Value = '64.739.51'
$ws.Range('E16').Value = '  +4.84%  '
$ws.Range('D17').Value = '3.474.13'
$ws.Range('E17').Value = '  +2.82%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '6.90'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('E18').Value = '  +3.80%  '
$ws.Range('D19').Value = '2.976.00'
$ws.Range('E19').Value = '  +1.99%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '446.34'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('E20').Value = '  +2.33%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.61'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('E21').Value = '  +2.36%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.676'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('E22').Value = '  +3.01%  '
$ws.Range('E23').Value = '  +5.49%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '80.90'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('E24').Value = '  +1.68%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '10.77'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('E25').Value = '  +5.95%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.24'
$ws.Range('D26').NumberFormat = 'General'
$ws.Range('E26').Value = '  +3.42%  '
$ws.Range('E27').Value = '  +7.02%  '
$ws.Range('E28').Value = '  -0.04%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.35'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('E29').Value = '  +14.17%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.66'
$ws.Range('D30').NumberFormat = 'General'
$ws.Range('E30').Value = '  +8.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0000105'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('E31').Value = '  +0.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.56'
$ws.Range('D32').NumberFormat = 'General'
$ws.Range('E32').Value = '  +2.18%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.109'
$ws.Range('D33').NumberFormat = 'General'
$ws.Range('E33').Value = '  +2.63%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '26.54'
$ws.Range('D34').NumberFormat = 'General'
$ws.Range('E34').Value = '  +3.59%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.980'
$ws.Range('D36').NumberFormat = 'General'
$ws.Range('E36').Value = '  +2.06%  '
$ws.Range('E37').Value = '  +3.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.10'
$ws.Range('D38').NumberFormat = 'General'
$ws.Range('E38').Value = '  +7.90%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '2.95'
$ws.Range('D39').NumberFormat = 'General'
$ws.Range('E39').Value = '  +4.13%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '48.82'
$ws.Range('D40').NumberFormat = 'General'
$ws.Range('E40').Value = '  -0.35%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '43.64'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('E41').Value = '  +10.51%  '
$ws.Range('E42').Value = '  +3.93%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.293'
$ws.Range('D43').NumberFormat = 'General'
$ws.Range('E43').Value = '  +8.82%  '
$ws.Range('E44').Value = '  +1.11%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '381.44'
$ws.Range('D45').NumberFormat = 'General'
$ws.Range('E45').Value = '  +13.02%  '
$ws.Range('D46').Value = '2.782.21'
$ws.Range('E46').Value = '  +3.09%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0346'
$ws.Range('D47').NumberFormat = 'General'
$ws.Range('E47').Value = '  +3.71%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '134.42'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('E48').Value = '  +0.48%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('B50').Value = 'Stellar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.105'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('E50').Value = '  +1.76%  '
$ws.Range('B51').Value = 'FLOKI'
$ws.Range('C51').Value = 'https://coinranking.com/coin/fmHk13Rqw+floki-floki'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.000219'
$ws.Range('D51').NumberFormat = 'General'
$ws.Range('E51').Value = '  +13.37%  '
